$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PSSM matrix (B2:K21) with supplemental-figure values
$data = New-Object 'object[,]' 20,10

$data[0,0] = -18.7395498972184
$data[0,1] = -0.716851576687688
$data[0,2] = -18.7395498972184
$data[0,3] = -18.7395498972184
$data[0,4] = -18.7395498972184
$data[0,5] = -18.7395498972184
$data[0,6] = -18.7395498972184
$data[0,7] = -18.7395498972184
$data[0,8] = -18.7395498972184
$data[0,9] = -18.7395498972184
$data[1,0] = -18.7395498972184
$data[1,1] = -18.7395498972184
$data[1,2] = -18.7395498972184
$data[1,3] = -18.7395498972184
$data[1,4] = -18.7395498972184
$data[1,5] = -18.7395498972184
$data[1,6] = -18.7395498972184
$data[1,7] = 0.431492482197306
$data[1,8] = -18.7395498972184
$data[1,9] = -18.7395498972184
$data[2,0] = -18.7395498972184
$data[2,1] = -0.6193474638810311
$data[2,2] = -0.05293652711367982
$data[2,3] = -18.7395498972184
$data[2,4] = 3.981168016705975
$data[2,5] = -18.7395498972184
$data[2,6] = 2.247464936394079
$data[2,7] = -18.7395498972184
$data[2,8] = 2.816456975166489
$data[2,9] = -18.7395498972184
$data[3,0] = -18.7395498972184
$data[3,1] = -0.1125540173803852
$data[3,2] = -18.7395498972184
$data[3,3] = -18.7395498972184
$data[3,4] = -18.7395498972184
$data[3,5] = 3.528542298685419
$data[3,6] = -18.7395498972184
$data[3,7] = -18.7395498972184
$data[3,8] = -18.7395498972184
$data[3,9] = -18.7395498972184
$data[4,0] = -18.7395498972184
$data[4,1] = -18.7395498972184
$data[4,2] = -18.7395498972184
$data[4,3] = -18.7395498972184
$data[4,4] = -18.7395498972184
$data[4,5] = -18.7395498972184
$data[4,6] = -18.7395498972184
$data[4,7] = -18.7395498972184
$data[4,8] = -18.7395498972184
$data[4,9] = -18.7395498972184
$data[5,0] = 4.321924963534633
$data[5,1] = -18.7395498972184
$data[5,2] = -18.7395498972184
$data[5,3] = -18.7395498972184
$data[5,4] = -18.7395498972184
$data[5,5] = -18.7395498972184
$data[5,6] = -18.7395498972184
$data[5,7] = -18.7395498972184
$data[5,8] = -18.7395498972184
$data[5,9] = -18.7395498972184
$data[6,0] = -18.7395498972184
$data[6,1] = -18.7395498972184
$data[6,2] = -18.7395498972184
$data[6,3] = 0.8585407785921646
$data[6,4] = -18.7395498972184
$data[6,5] = -18.7395498972184
$data[6,6] = -18.7395498972184
$data[6,7] = -18.7395498972184
$data[6,8] = -18.7395498972184
$data[6,9] = -18.7395498972184
$data[7,0] = -18.7395498972184
$data[7,1] = -18.7395498972184
$data[7,2] = -18.7395498972184
$data[7,3] = -18.7395498972184
$data[7,4] = -18.7395498972184
$data[7,5] = -18.7395498972184
$data[7,6] = -18.7395498972184
$data[7,7] = -18.7395498972184
$data[7,8] = -18.7395498972184
$data[7,9] = -18.7395498972184
$data[8,0] = -18.7395498972184
$data[8,1] = -18.7395498972184
$data[8,2] = -18.7395498972184
$data[8,3] = -18.7395498972184
$data[8,4] = -18.7395498972184
$data[8,5] = -18.7395498972184
$data[8,6] = -18.7395498972184
$data[8,7] = 0.2686078658242152
$data[8,8] = -18.7395498972184
$data[8,9] = 2.023799460236277
$data[9,0] = -18.7395498972184
$data[9,1] = -18.7395498972184
$data[9,2] = -18.7395498972184
$data[9,3] = 2.768918545928209
$data[9,4] = -18.7395498972184
$data[9,5] = 1.450044252477357
$data[9,6] = -18.7395498972184
$data[9,7] = -18.7395498972184
$data[9,8] = -18.7395498972184
$data[9,9] = 1.273332390812999
$data[10,0] = -18.7395498972184
$data[10,1] = -18.7395498972184
$data[10,2] = -18.7395498972184
$data[10,3] = -18.7395498972184
$data[10,4] = -18.7395498972184
$data[10,5] = -18.7395498972184
$data[10,6] = -18.7395498972184
$data[10,7] = -18.7395498972184
$data[10,8] = -18.7395498972184
$data[10,9] = -18.7395498972184
$data[11,0] = -18.7395498972184
$data[11,1] = -18.7395498972184
$data[11,2] = -18.7395498972184
$data[11,3] = 1.74848340341479
$data[11,4] = -18.7395498972184
$data[11,5] = -18.7395498972184
$data[11,6] = -18.7395498972184
$data[11,7] = -18.7395498972184
$data[11,8] = 0.9804889018983168
$data[11,9] = 2.817560839567585
$data[12,0] = -18.7395498972184
$data[12,1] = -18.7395498972184
$data[12,2] = 0.2671295571050032
$data[12,3] = -18.7395498972184
$data[12,4] = -18.7395498972184
$data[12,5] = -18.7395498972184
$data[12,6] = -18.7395498972184
$data[12,7] = -18.7395498972184
$data[12,8] = -18.7395498972184
$data[12,9] = 1.745498125290223
$data[13,0] = -18.7395498972184
$data[13,1] = -18.7395498972184
$data[13,2] = -1.017936658376097
$data[13,3] = -18.7395498972184
$data[13,4] = -18.7395498972184
$data[13,5] = -18.7395498972184
$data[13,6] = -18.7395498972184
$data[13,7] = -18.7395498972184
$data[13,8] = -18.7395498972184
$data[13,9] = -18.7395498972184
$data[14,0] = -18.7395498972184
$data[14,1] = -18.7395498972184
$data[14,2] = -18.7395498972184
$data[14,3] = -18.7395498972184
$data[14,4] = -18.7395498972184
$data[14,5] = -18.7395498972184
$data[14,6] = -18.7395498972184
$data[14,7] = -18.7395498972184
$data[14,8] = 1.961774334312111
$data[14,9] = -18.7395498972184
$data[15,0] = -18.7395498972184
$data[15,1] = 0.1688691123605907
$data[15,2] = -0.5286835821120534
$data[15,3] = -18.7395498972184
$data[15,4] = -18.7395498972184
$data[15,5] = -18.7395498972184
$data[15,6] = 2.073840667755666
$data[15,7] = 0.1411606156328546
$data[15,8] = 1.833953991366581
$data[15,9] = -18.7395498972184
$data[16,0] = -18.7395498972184
$data[16,1] = -18.7395498972184
$data[16,2] = -18.7395498972184
$data[16,3] = -18.7395498972184
$data[16,4] = -18.7395498972184
$data[16,5] = -18.7395498972184
$data[16,6] = 2.067783797110129
$data[16,7] = -0.413045798343573
$data[16,8] = 1.816391933078639
$data[16,9] = -18.7395498972184
$data[17,0] = -18.7395498972184
$data[17,1] = -18.7395498972184
$data[17,2] = 2.882328930494333
$data[17,3] = -18.7395498972184
$data[17,4] = -18.7395498972184
$data[17,5] = -18.7395498972184
$data[17,6] = 1.585941127414729
$data[17,7] = 1.266552879878486
$data[17,8] = -18.7395498972184
$data[17,9] = -18.7395498972184
$data[18,0] = -18.7395498972184
$data[18,1] = 3.58472103700117
$data[18,2] = 3.212908763996816
$data[18,3] = -18.7395498972184
$data[18,4] = 2.07294785806998
$data[18,5] = -18.7395498972184
$data[18,6] = 0.9748985522544372
$data[18,7] = 3.721045230698457
$data[18,8] = -18.7395498972184
$data[18,9] = 1.638508890801926
$data[19,0] = -18.7395498972184
$data[19,1] = 2.230652209934514
$data[19,2] = -18.7395498972184
$data[19,3] = 3.001935605326898
$data[19,4] = -18.7395498972184
$data[19,5] = 2.518038789244538
$data[19,6] = 0.9118609790755273
$data[19,7] = -18.7395498972184
$data[19,8] = -18.7395498972184
$data[19,9] = -18.7395498972184

$ws.Range("B2:K21").Value = $data

Write-Host "Updated PSSM values in B2:K21"